# Replace the three-digit x one-digit multiplication problems/answers
# in the table with the new values described by the diff.

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "653×3=1959"; New = "530×8=4240" },
    @{ Old = "579×8=4632"; New = "743×9=6687" },
    @{ Old = "264×9=2376"; New = "797×2=1594" },
    @{ Old = "910×6=5460"; New = "877×8=7016" },
    @{ Old = "211×9=1899"; New = "209×3=627" },
    @{ Old = "459×9=4131"; New = "229×9=2061" },
    @{ Old = "779×6=4674"; New = "382×9=3438" },
    @{ Old = "760×8=6080"; New = "902×6=5412" },
    @{ Old = "983×9=8847"; New = "339×3=1017" },
    @{ Old = "445×7=3115"; New = "597×5=2985" },
    @{ Old = "386×7=2702"; New = "145×2=290" },
    @{ Old = "971×9=8739"; New = "528×8=4224" },
    @{ Old = "131×9=1179"; New = "139×2=278" },
    @{ Old = "499×7=3493"; New = "549×6=3294" },
    @{ Old = "183×8=1464"; New = "161×7=1127" },
    @{ Old = "290×5=1450"; New = "551×3=1653" },
    @{ Old = "496×7=3472"; New = "140×3=420" },
    @{ Old = "200×2=400";  New = "657×8=5256" },
    @{ Old = "327×3=981";  New = "811×4=3244" },
    @{ Old = "848×2=1696"; New = "163×5=815" },
    @{ Old = "750×7=5250"; New = "855×6=5130" },
    @{ Old = "438×5=2190"; New = "154×3=462" },
    @{ Old = "886×2=1772"; New = "881×6=5286" },
    @{ Old = "126×6=756";  New = "181×9=1629" },
    @{ Old = "905×4=3620"; New = "309×2=618" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $true, $false, $false, $false, $true, 1, $false, $r.New, 2)
}

$d.Save()
